$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: NIK, NAMA, KODE_LOKASI (replacing NAMA, PROVINSI, KABUPATEN, KECAMATAN, KELURAHAN)
$ws.Range("A1").Value = "NIK"
$ws.Range("B1").Value = "NAMA"
$ws.Range("C1").Value = "KODE_LOKASI"

# Remove the now-unused D1/E1 header cells entirely
$ws.Range("D1:E1").Clear()

# Move the active selection to C12, matching the saved cursor position
$ws.Range("C12").Select() | Out-Null
